$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Start/Finish (columns C and D) date values for the task rows,
# while leaving their (date) cell formatting intact.
$rows = @(3, 4, 5, 8, 9, 10, 13, 14, 15, 18, 19, 20)
foreach ($r in $rows) {
    $ws.Range("C$r`:D$r").Value = $null
}

# Update the saved selection to match the author's edit.
$ws.Range("C3:D20").Select() | Out-Null
